$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-format cells keep the "@" number format (style index 1) like
# their header/sibling cells, matching cells that pre-existed as s="1".
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "127.0.0.1"

$ws.Range("A2").Value = "GameServer_1"

$ws.Range("B2").Value = "000104001"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "GameServer_1"

$ws.Range("D2").Value = 5000
$ws.Range("E2").Value = 1
$ws.Range("G2").Value = 4001

# Update the data validation range to start from F3 (row 2 now has a concrete value)
$ws.Range("F2:F1048576").Validation.Delete()
$ws.Range("F3:F1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# Update the active selection to G3
$ws.Range("G3").Select()
